# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.541.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.724.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3725"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.31%  "

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.393"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.724.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001073"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06651"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.149"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.519.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.92%  "

$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.410"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +21.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.395"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.918.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.105"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.970"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08627"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.76%  "

$ws.Range("E35").Value = "  +3.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.365"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02332"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.417"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.222"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6204"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6016"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.51%  "
